# Research Project 20% - apply "strike-through completed requirements" +
# split trailing word into its own run (grammar-check style) + misc cleanup.

# --- 1. "It must be possible to create and delete objects" -> strikethrough
$d = $word.ActiveDocument
$d.Paragraphs(7).Range.Font.StrikeThrough = 1
Write-Host "strike p7 done"

# --- 2. "Objects must be serialized ..." -> strikethrough, drop _GoBack bookmark
$d = $word.ActiveDocument
$d.Paragraphs(8).Range.Font.StrikeThrough = 1
Write-Host "strike p8 done"

$d = $word.ActiveDocument
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
Write-Host "bookmark removed"

# --- 3. Split trailing word of the GUI bullet points into their own run
#        (mirrors Word's "flag last word" proofing-split pattern).
$d = $word.ActiveDocument
$p = $d.Paragraphs(11)
$s = $p.Range.Start
$r = $d.Range($s + 12, $s + 19)
Write-Host "p11 word: $($r.Text)"
$r.Font.Size = 16
$d = $word.ActiveDocument
$r = $d.Range($s + 12, $s + 19)
$r.Font.Size = 14
Write-Host "p11 split done"

$d = $word.ActiveDocument
$p = $d.Paragraphs(12)
$s = $p.Range.Start
$r = $d.Range($s + 15, $s + 21)
Write-Host "p12 word: $($r.Text)"
$r.Font.Size = 16
$d = $word.ActiveDocument
$r = $d.Range($s + 15, $s + 21)
$r.Font.Size = 14
Write-Host "p12 split done"

$d = $word.ActiveDocument
$p = $d.Paragraphs(13)
$s = $p.Range.Start
$r = $d.Range($s + 18, $s + 25)
Write-Host "p13 word: $($r.Text)"
$r.Font.Size = 16
$d = $word.ActiveDocument
$r = $d.Range($s + 18, $s + 25)
$r.Font.Size = 14
Write-Host "p13 split done"

$d = $word.ActiveDocument
$p = $d.Paragraphs(14)
$s = $p.Range.Start
$r = $d.Range($s + 16, $s + 23)
Write-Host "p14 word: $($r.Text)"
$r.Font.Size = 16
$d = $word.ActiveDocument
$r = $d.Range($s + 16, $s + 23)
$r.Font.Size = 14
Write-Host "p14 split done"

# --- 4. Merge the "Don't" / " forget to comment code!" runs into a single run
$d = $word.ActiveDocument
$d.Content.Find.Execute("Don" + [char]0x2019 + "t forget to comment code!", $true, $false, $false, $false, $false, $true, 1, $false, "Don" + [char]0x2019 + "t forget to comment code!", 2)
Write-Host "merge done"

# --- 5. Split trailing word of the final "Screencast..." line into its own run
$d = $word.ActiveDocument
$p = $d.Paragraphs(19)
$s = $p.Range.Start
$r = $d.Range($s + 26, $s + 35)
Write-Host "p19 word: $($r.Text)"
$r.Font.Size = 16
$d = $word.ActiveDocument
$r = $d.Range($s + 26, $s + 35)
$r.Font.Size = 14
Write-Host "p19 split done"

Write-Host "done"
